# oppdatert med siste forbrukstall
# Appends the daily kWh consumption readings for 20.12.2023 - 03.01.2024
# to the bottom of the existing "dato / kWh" table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (dato, kWh) -------------------------------------------
# Rows 51-62 -> 20.12.2023 .. 31.12.2023 (pasted in one batch)
$data1 = @(
    @("20.12.2023", 904.2),
    @("21.12.2023", 1038.9000000000001),
    @("22.12.2023", 960.1),
    @("23.12.2023", 1296.8499999999999),
    @("24.12.2023", 1334.1),
    @("25.12.2023", 870.8),
    @("26.12.2023", 883.7),
    @("27.12.2023", 984.85),
    @("28.12.2023", 977.25),
    @("29.12.2023", 868.9),
    @("30.12.2023", 926.2),
    @("31.12.2023", 815.95)
)

# Rows 63-65 -> 01.01.2024 .. 03.01.2024 (pasted in a second batch)
$data2 = @(
    @("01.01.2024", 780.8),
    @("02.01.2024", 856.45),
    @("03.01.2024", 1393.2)
)

$row = 51
foreach ($pair in $data1) {
    # Leading apostrophe forces literal text so ambiguous dd.mm.yyyy values
    # (day <= 12) aren't silently reinterpreted as dates.
    $ws.Cells.Item($row, 1).Value = "'" + $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
foreach ($pair in $data2) {
    $ws.Cells.Item($row, 1).Value = "'" + $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# --- View state: scrolled down to show the newly added rows --------------
$ws.Range("E45").Select()
$excel.ActiveWindow.ScrollRow = 26
